$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# 1. Title shape: "A" + " " + "slide" -> "A slide"
$s.Shapes.Item("Title 1").TextFrame.TextRange.Text = "A slide"

# 2. Table cell: "a" + " " + "table" -> "a table"
$tblShape = $s.Shapes.Item("Content Placeholder 5")
$tbl = $tblShape.Table
$tbl.Cell(1, 2).Shape.TextFrame.TextRange.Text = "a table"

# 3. TextBox: "Plus" + " " + "an" + " " + "image" -> "Plus an image"
$s.Shapes.Item("TextBox 3").TextFrame.TextRange.Text = "Plus an image"
